$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (G:K) before the existing "arrecadado_sucesso" block.
# This shifts every existing column from G onward five places to the right
# and (per this engine's Insert behaviour) copies the style of column G
# ("R$ #,##0.00" / style 3) onto the data rows of the freshly inserted cells,
# and the header style (style 4) onto row 1.
$ws.Range("G:K").Insert()

# New header row values for the inserted columns.
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# New data values (row 2 = sub / apoia.se).
$ws.Range("G2").Value = 165199.0578149446
$ws.Range("H2").Value = 1205.832538795216
$ws.Range("I2").Value = 2163.288658625353
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 21176.91783511972

# New data values (row 3 = sub / catarse).
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Make sure the new data cells carry the same "R$ #,##0.00" currency format
# used by the rest of the monetary columns (matches style index 3).
$ws.Range("G2:K3").NumberFormat = "R$ #,##0.00"
